# Remove column M from the alcohol measurement data sheet.
# The old column M ("M") is deleted entirely; the old column N shifts left
# to become the new column M, and the used range shrinks from A1:N119 to
# A1:M119.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Select column M (13th column) first so the active cell / selection after
# the delete lands on M1, mirroring how Excel leaves the selection after an
# interactive "Delete Column" on a selected column.
$ws.Columns.Item(13).Select()
$ws.Columns.Item(13).Delete()

# After deletion, make sure the active selection is the single cell M1.
$ws.Range("M1").Select()
